# ADD results from server
# Update investment cost result cells on the "2025", "2030" and "2035" sheets
# with refreshed values received from the server.

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 2778.902526399997
$ws.Range("E2").Value = 290927.2506141524
$ws.Range("G2").Value = 80959.25712662093
$ws.Range("I2").Value = 148652.5872276
$ws.Range("L2").Value = 509125.9821312751
$ws.Range("M2").Value = 112470.9127927
$ws.Range("N2").Value = 71977.22211759936
$ws.Range("O2").Value = 68708.80120585456

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 52284.38913067689
$ws.Range("E2").Value = 269123.7921032493
$ws.Range("I2").Value = 219942.2466224718
$ws.Range("L2").Value = 229245.3249170133
$ws.Range("M2").Value = 105590.95466293
$ws.Range("N2").Value = 36433.63837700079
$ws.Range("O2").Value = 25074.86257357796

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 22324.02753895484
$ws.Range("B2").Value = 14448.48472786732
$ws.Range("E2").Value = 112847.7913764417
$ws.Range("I2").Value = 163473.2120365721
$ws.Range("M2").Value = 57872.79574411505
$ws.Range("N2").Value = 49753.5561971986
$ws.Range("O2").Value = 58474.09285693887
